$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date column U (2025-11-22) -----------------------------------
# Force the new header cell to hold the literal text "2025-11-22" rather
# than letting Excel auto-convert it to a date serial, then copy the
# existing header formatting (bold/border/centered) from T1 onto it.
$ws.Range("U1").NumberFormat = "@"
$ws.Range("U1").Value = "2025-11-22"
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)

# New column U attendance marks - everyone absent on this date
$ws.Range("U2").Value = "❌"
$ws.Range("U3").Value = "❌"
$ws.Range("U4").Value = "❌"
$ws.Range("U5").Value = "❌"
$ws.Range("U6").Value = "❌"
$ws.Range("U7").Value = "❌"
$ws.Range("U8").Value = "❌"
$ws.Range("U9").Value = "❌"

# --- Updated Total (S) day-counts now that there's one more date ------
$ws.Range("S2").Value = 16
$ws.Range("S3").Value = 16
$ws.Range("S4").Value = 16
$ws.Range("S5").Value = 16
$ws.Range("S6").Value = 16
$ws.Range("S7").Value = 16
$ws.Range("S8").Value = 1
$ws.Range("S9").Value = 16

# --- Updated Attendance % (T) recalculated against the new totals -----
$ws.Range("T6").Value = 6.2
$ws.Range("T9").Value = 6.2
